$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 105698
$ws.Range("B2").Value = "B101"

# Row 3
$ws.Range("A3").Value = 105698
$ws.Range("B3").Value = "S354"

# Row 4
$ws.Range("A4").Value = 105698
$ws.Range("B4").Value = "S489"
$ws.Range("C4").Value = 13

# Row 5
$ws.Range("A5").Value = 104891
$ws.Range("B5").Value = "S354"

# Row 6
$ws.Range("A6").Value = 104891
$ws.Range("B6").Value = "S489"

# Update selection to D8 (as reflected in the sheetView selection)
$ws.Range("D8").Select() | Out-Null
